# Applies the "Updated symbol list" GitHub Actions refresh to the
# cryptocurrency tracking sheet: refreshed prices/volumes, a handful
# of tokens re-ranked (rows 7-20 shuffle up/down a slot), and the
# snapshot hour bumped from 18 to 19 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: a row number plus the column letters on that row whose
# text needs to change (Coin / Link / Price / Volume(1h) / Hora).
$rowUpdates = @(
    @{ Row = 2; Cells = @{ "D" = "246.56"; "E" = "0.67%"; "G" = "19" } }
    @{ Row = 3; Cells = @{ "D" = "29.93"; "E" = "10.27%"; "G" = "19" } }
    @{ Row = 4; Cells = @{ "D" = "5.175"; "E" = "1.68%"; "G" = "19" } }
    @{ Row = 5; Cells = @{ "D" = "0.05728"; "G" = "19" } }
    @{ Row = 6; Cells = @{ "D" = "6.584"; "E" = "1.13%"; "G" = "19" } }
    @{ Row = 7; Cells = @{ "B" = "MXToken"; "C" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; "D" = "0.8560"; "E" = "4.44%"; "G" = "19" } }
    @{ Row = 8; Cells = @{ "B" = "FTXToken"; "C" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; "D" = "0.8672"; "E" = "0.43%"; "G" = "19" } }
    @{ Row = 9; Cells = @{ "B" = "WazirX"; "C" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; "D" = "0.1361"; "E" = "2.24%"; "G" = "19" } }
    @{ Row = 10; Cells = @{ "B" = "MandalaExchangeToken"; "C" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; "D" = "0.07062"; "E" = "1.81%"; "G" = "19" } }
    @{ Row = 11; Cells = @{ "B" = "BitrueCoin"; "C" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; "D" = "0.02926"; "E" = "3.59%"; "G" = "19" } }
    @{ Row = 12; Cells = @{ "B" = "BitMartToken"; "C" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; "D" = "0.09387"; "E" = "-0.14%"; "G" = "19" } }
    @{ Row = 13; Cells = @{ "B" = "BitForexToken"; "C" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; "D" = "0.001519"; "E" = "-0.04%"; "G" = "19" } }
    @{ Row = 14; Cells = @{ "B" = "CoinExToken"; "C" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; "D" = "0.04162"; "E" = "2.37%"; "G" = "19" } }
    @{ Row = 15; Cells = @{ "B" = "One"; "C" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; "D" = "0.0006028"; "E" = "0.51%"; "G" = "19" } }
    @{ Row = 16; Cells = @{ "D" = "0.006066"; "E" = "-1.71%"; "G" = "19" } }
    @{ Row = 17; Cells = @{ "E" = "5,070.51%"; "G" = "19" } }
    @{ Row = 18; Cells = @{ "D" = "3.488"; "E" = "-0.54%"; "G" = "19" } }
    @{ Row = 19; Cells = @{ "B" = "GateToken"; "C" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; "D" = "3.096"; "E" = "2.85%"; "G" = "19" } }
    @{ Row = 20; Cells = @{ "B" = "BTSEToken"; "C" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; "D" = "2.267"; "E" = "1.72%"; "G" = "19" } }
    @{ Row = 21; Cells = @{ "D" = "0.3155"; "E" = "-0.69%"; "G" = "19" } }
    @{ Row = 22; Cells = @{ "E" = "6.43%"; "G" = "19" } }
    @{ Row = 23; Cells = @{ "D" = "0.1317"; "E" = "1.05%"; "G" = "19" } }
    @{ Row = 24; Cells = @{ "D" = "3.467"; "E" = "-2.57%"; "G" = "19" } }
    @{ Row = 25; Cells = @{ "D" = "0.1380"; "E" = "0.45%"; "G" = "19" } }
    @{ Row = 26; Cells = @{ "D" = "0.005017"; "E" = "26.19%"; "G" = "19" } }
    @{ Row = 27; Cells = @{ "D" = "0.001222"; "E" = "0.26%"; "G" = "19" } }
    @{ Row = 28; Cells = @{ "D" = "0.0001210"; "E" = "22.27%"; "G" = "19" } }
    @{ Row = 29; Cells = @{ "G" = "19" } }
    @{ Row = 30; Cells = @{ "G" = "19" } }
    @{ Row = 31; Cells = @{ "G" = "19" } }
    @{ Row = 32; Cells = @{ "G" = "19" } }
    @{ Row = 33; Cells = @{ "G" = "19" } }
    @{ Row = 34; Cells = @{ "G" = "19" } }
    @{ Row = 35; Cells = @{ "G" = "19" } }
    @{ Row = 36; Cells = @{ "G" = "19" } }
    @{ Row = 37; Cells = @{ "G" = "19" } }
    @{ Row = 38; Cells = @{ "G" = "19" } }
    @{ Row = 39; Cells = @{ "G" = "19" } }
    @{ Row = 40; Cells = @{ "D" = "0.03753"; "E" = "0.80%"; "G" = "19" } }
    @{ Row = 41; Cells = @{ "D" = "0.005749"; "E" = "67.38%"; "G" = "19" } }
    @{ Row = 42; Cells = @{ "D" = "0.1073"; "E" = "1.39%"; "G" = "19" } }
    @{ Row = 43; Cells = @{ "D" = "0.002000"; "E" = "-19.33%"; "G" = "19" } }
    @{ Row = 44; Cells = @{ "D" = "0.009666"; "E" = "3.11%"; "G" = "19" } }
    @{ Row = 45; Cells = @{ "D" = "0.00005231"; "E" = "1.65%"; "G" = "19" } }
    @{ Row = 46; Cells = @{ "E" = "0.01%"; "G" = "19" } }
    @{ Row = 47; Cells = @{ "D" = "0.06469"; "E" = "-45.07%"; "G" = "19" } }
    @{ Row = 48; Cells = @{ "D" = "0.002519"; "E" = "-0.40%"; "G" = "19" } }
    @{ Row = 49; Cells = @{ "E" = "0.01%"; "G" = "19" } }
    @{ Row = 50; Cells = @{ "D" = "0.0002000"; "E" = "0.01%"; "G" = "19" } }
    @{ Row = 51; Cells = @{ "G" = "19" } }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $cellRef = "$col$row"
        # Every touched column in this sheet stores plain text (prices,
        # percentages, names and URLs alike), so force text format first
        # -- otherwise Excel would auto-coerce values like "0.67%" or
        # "246.56" into numeric/percentage cells instead of literal text.
        $ws.Range($cellRef).NumberFormat = "@"
        $ws.Range($cellRef).Value = $update.Cells[$col]
    }
}
